# Add the new "Lemon VS Orange" worksheet as the first sheet in the workbook,
# matching the structural shift seen in the diff: the new sheet becomes
# sheet1.xml / rId1, and Dataset / Descriptive stats / Correlation shift
# down by one position (rId2 / rId3 / rId4 respectively).
$wb = $excel.ActiveWorkbook

$ds = $wb.Worksheets.Item("Dataset")
$newSheet = $wb.Worksheets.Add($ds)
$newSheet.Name = "Lemon VS Orange"

# --- formatting: reuse the existing "header row" (italic, centered, top+bottom
# border) and "final row" (bottom border) styles already present in the
# workbook (used on the Descriptive stats / Correlation sheets) by copying
# them onto our new ranges before writing values, so no redundant style
# entries get created. ---
$statsSheet = $wb.Worksheets.Item("Descriptive stats")
$statsSheet.Range("A1:C1").Copy($newSheet.Range("A3:C3"))
$statsSheet.Range("A15:C15").Copy($newSheet.Range("A14:C14"))

# --- content: t-Test: Two-Sample Assuming Equal Variances, comparing the
# Lemon and Orange sales columns from the Dataset sheet ---
$newSheet.Range("A1").Value = "Test d'égalité des espérances: deux observations de variances égales"

$newSheet.Range("A3").Value = ""
$newSheet.Range("B3").Value = "Lemon"
$newSheet.Range("C3").Value = "Orange"

$newSheet.Range("A4").Value = "Moyenne"
$newSheet.Range("B4").Value = 116.58064516129032
$newSheet.Range("C4").Value = 80.354838709677423

$newSheet.Range("A5").Value = "Variance"
$newSheet.Range("B5").Value = 683.11827956989293
$newSheet.Range("C5").Value = 489.7698924731182

$newSheet.Range("A6").Value = "Observations"
$newSheet.Range("B6").Value = 31
$newSheet.Range("C6").Value = 31

$newSheet.Range("A7").Value = "Variance pondérée"
$newSheet.Range("B7").Value = 586.44408602150554

$newSheet.Range("A8").Value = "Différence hypothétique des moyennes"
$newSheet.Range("B8").Value = 0

$newSheet.Range("A9").Value = "Degré de liberté"
$newSheet.Range("B9").Value = 60

$newSheet.Range("A10").Value = "Statistique t"
$newSheet.Range("B10").Value = 5.8893939518238767

$newSheet.Range("A11").Value = "P(T<=t) unilatéral"
$newSheet.Range("B11").Value = 0.0000000939311262965143682

$newSheet.Range("A12").Value = "Valeur critique de t (unilatéral)"
$newSheet.Range("B12").Value = 1.6706488649046354

$newSheet.Range("A13").Value = "P(T<=t) bilatéral"
$newSheet.Range("B13").Value = 0.0000001878622525930287363

$newSheet.Range("A14").Value = "Valeur critique de t (bilatéral)"
$newSheet.Range("B14").Value = 2.0002978220142609
$newSheet.Range("C14").Value = ""

# --- column widths (approximate Excel's best-fit character widths) ---
$newSheet.Columns.Item(1).ColumnWidth = 58.45
$newSheet.Columns.Item(2).ColumnWidth = 13.63
$newSheet.Columns.Item(3).ColumnWidth = 11.82

# --- Dataset sheet view: selection moves from E1:H32 to C3 ---
$ds.Activate()
$ds.Range("C3").Select()

# --- restore the new sheet as the active / selected sheet, with the whole
# populated range selected (A1:C14) ---
$newSheet.Activate()
$newSheet.Range("A1:C14").Select()
